$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1882.3636
$ws.Range("J19").Value = 1995.6
$ws.Range("L19").Value = 1995.6
$ws.Range("N19").Value = -2345.6
$ws.Range("H107").Value = 447.9
$ws.Range("I107").Value = 447.9
$ws.Range("K107").Value = 447.9
$ws.Range("M107").Value = 1472.1
$ws.Range("H112").Value = 5116.4863
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").ClearContents()
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20082.217
$ws.Range("I32").Value = 19252.334
$ws.Range("K32").Value = 19252.334
$ws.Range("M32").Value = -18965.334
$ws.Range("H61").Value = 6206.4287
$ws.Range("I61").Value = 5754.4736
$ws.Range("K61").Value = 5754.4736
$ws.Range("M61").Value = -5542.4736
$ws.Range("H122").Value = 2901.3489
$ws.Range("I122").Value = 2208.7632
$ws.Range("J122").Value = 8165
$ws.Range("K122").Value = 6626.2896
$ws.Range("L122").Value = 24495
$ws.Range("M122").Value = -4176.2896
$ws.Range("N122").Value = -29395
$ws.Range("H132").Value = 4754.2144
$ws.Range("I132").Value = 3796.7917
$ws.Range("K132").Value = 11390.3751
$ws.Range("M132").Value = -8860.375100000001
$ws.Range("H136").Value = 6206.4287
$ws.Range("I136").Value = 5754.4736
$ws.Range("K136").Value = 17263.4208
$ws.Range("M136").Value = -14713.4208
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1117.1578
$ws.Range("I94").Value = 587.4286
$ws.Range("J94").Value = 2600.4
$ws.Range("K94").Value = 587.4286
$ws.Range("L94").Value = 2600.4
$ws.Range("M94").Value = -136.4286
$ws.Range("N94").Value = -3502.4
$ws.Range("H134").Value = 5468.3257
$ws.Range("I134").Value = 4188.359
$ws.Range("J134").Value = 17948
$ws.Range("K134").Value = 12565.077
$ws.Range("L134").Value = 53844
$ws.Range("M134").Value = -10030.077
$ws.Range("N134").Value = -58914
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9468.875
$ws.Range("I31").Value = 3628.7368
$ws.Range("K31").Value = 3628.7368
$ws.Range("M31").Value = -3333.7368
$ws.Range("H34").Value = 9468.875
$ws.Range("I34").Value = 3628.7368
$ws.Range("K34").Value = 3628.7368
$ws.Range("M34").Value = -3426.7368
$ws.Range("H51").Value = 53999.5
$ws.Range("J51").Value = 53999.5
$ws.Range("L51").Value = 53999.5
$ws.Range("N51").Value = -55471.5
$ws.Range("H52").Value = 70565.75
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 70565.75
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 70565.75
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -71153.75
$ws.Range("H59").Value = 78584.86
$ws.Range("J59").Value = 101998
$ws.Range("L59").Value = 101998
$ws.Range("N59").Value = -104288
$ws.Range("H61").Value = 53999.5
$ws.Range("J61").Value = 53999.5
$ws.Range("L61").Value = 53999.5
$ws.Range("N61").Value = -54695.5
$ws.Range("H95").Value = 32579.5
$ws.Range("J95").Value = 32579.5
$ws.Range("L95").Value = 32579.5
$ws.Range("N95").Value = -38071.5
$ws.Range("H99").Value = 5787.896
$ws.Range("I99").Value = 6890.32
$ws.Range("J99").Value = 4589.609
$ws.Range("K99").Value = 6890.32
$ws.Range("L99").Value = 4589.609
$ws.Range("M99").Value = -5392.32
$ws.Range("N99").Value = -7585.609
$ws.Range("H126").Value = 5787.896
$ws.Range("I126").Value = 6890.32
$ws.Range("J126").Value = 4589.609
$ws.Range("K126").Value = 20670.96
$ws.Range("L126").Value = 13768.827
$ws.Range("M126").Value = -18200.96
$ws.Range("N126").Value = -18708.827
$ws.Range("H132").Value = 5439.7085
$ws.Range("I132").Value = 3996
$ws.Range("K132").Value = 11988
$ws.Range("M132").Value = -9458
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 88
$ws.Range("I21").Value = 88
$ws.Range("K21").Value = 264
$ws.Range("M21").Value = -91
$ws.Range("H23").Value = 88
$ws.Range("J23").Value = 121
$ws.Range("L23").Value = 363
$ws.Range("N23").Value = -833
$ws.Range("H32").Value = 2921.158
$ws.Range("I32").Value = 1499
$ws.Range("J32").Value = 3000.1667
$ws.Range("K32").Value = 4497
$ws.Range("L32").Value = 9000.500100000001
$ws.Range("M32").Value = -4214
$ws.Range("N32").Value = -9566.500100000001
$ws.Range("H55").Value = 2625
$ws.Range("I55").Value = 500
$ws.Range("J55").Value = 3333.3333
$ws.Range("K55").Value = 1500
$ws.Range("L55").Value = 9999.999899999999
$ws.Range("M55").Value = -1323
$ws.Range("N55").Value = -10353.9999
$ws.Range("H58").Value = 861
$ws.Range("I58").Value = 912.5
$ws.Range("J58").Value = 655
$ws.Range("K58").Value = 2737.5
$ws.Range("L58").Value = 1965
$ws.Range("M58").Value = -2609.5
$ws.Range("N58").Value = -2221
$ws.Range("H107").Value = 497.55554
$ws.Range("J107").Value = 579.0769
$ws.Range("L107").Value = 1737.2307
$ws.Range("N107").Value = -5577.2307
$ws.Range("H132").Value = 2612.2856
$ws.Range("I132").Value = 2684
$ws.Range("J132").Value = 2558.5
$ws.Range("K132").Value = 24156
$ws.Range("L132").Value = 23026.5
$ws.Range("M132").Value = -21626
$ws.Range("N132").Value = -28086.5
$ws.Range("H137").Value = 4466.5
$ws.Range("J137").Value = 7405
$ws.Range("L137").Value = 22215
$ws.Range("N137").Value = -32415
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 11251.25
$ws.Range("I55").Value = 7001.6665
$ws.Range("K55").Value = 7001.6665
$ws.Range("M55").Value = -6674.6665
$ws.Range("H70").Value = 4999
$ws.Range("J70").Value = 4999
$ws.Range("L70").Value = 4999
$ws.Range("N70").Value = -5539
$ws.Range("H73").Value = 4999
$ws.Range("J73").Value = 4999
$ws.Range("L73").Value = 4999
$ws.Range("N73").Value = -6871
$ws.Range("H126").Value = 11628.36
$ws.Range("I126").Value = 11123.182
$ws.Range("K126").Value = 33369.546
$ws.Range("M126").Value = -30899.546
$ws.Range("H132").Value = 9016
$ws.Range("I132").Value = 6243.5557
$ws.Range("J132").Value = 17333.334
$ws.Range("K132").Value = 18730.6671
$ws.Range("L132").Value = 52000.00199999999
$ws.Range("M132").Value = -16200.6671
$ws.Range("N132").Value = -57060.00199999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 20021
$ws.Range("I34").Value = 20021
$ws.Range("K34").Value = 20021
$ws.Range("M34").Value = -19849
$ws.Range("H40").Value = 5185.2
$ws.Range("I40").Value = 3761.3333
$ws.Range("K40").Value = 3761.3333
$ws.Range("M40").Value = -3625.3333
$ws.Range("H82").Value = 2544
$ws.Range("I82").Value = 1026.2858
$ws.Range("J82").Value = 4314.6665
$ws.Range("K82").Value = 1026.2858
$ws.Range("L82").Value = 4314.6665
$ws.Range("M82").Value = -665.2858000000001
$ws.Range("N82").Value = -5036.6665
$ws.Range("H85").Value = 2544
$ws.Range("I85").Value = 1026.2858
$ws.Range("J85").Value = 4314.6665
$ws.Range("K85").Value = 1026.2858
$ws.Range("L85").Value = 4314.6665
$ws.Range("M85").Value = 221.7141999999999
$ws.Range("N85").Value = -6810.6665
$ws.Range("H93").Value = 3179.818
$ws.Range("I93").Value = 1733.75
$ws.Range("K93").Value = 1733.75
$ws.Range("M93").Value = -485.75
$ws.Range("H122").Value = 3603.2
$ws.Range("I122").Value = 2052.5454
$ws.Range("J122").Value = 7867.5
$ws.Range("K122").Value = 6157.6362
$ws.Range("L122").Value = 23602.5
$ws.Range("M122").Value = -3707.6362
$ws.Range("N122").Value = -28502.5
$ws.Range("H136").Value = 9178.361000000001
$ws.Range("I136").Value = 7746.4546
$ws.Range("K136").Value = 23239.3638
$ws.Range("M136").Value = -20689.3638
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 39999
$ws.Range("I34").Value = 39999
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 39999
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -39796
$ws.Range("N34").ClearContents()
$ws.Range("H122").Value = 3972.2
$ws.Range("I122").Value = 4024.6667
$ws.Range("K122").Value = 12074.0001
$ws.Range("M122").Value = -9624.000100000001
$ws.Range("H126").Value = 3768.8
$ws.Range("I126").Value = 1698
$ws.Range("J126").Value = 6875
$ws.Range("K126").Value = 5094
$ws.Range("L126").Value = 20625
$ws.Range("M126").Value = -2624
$ws.Range("N126").Value = -25565
$ws.Range("H132").Value = 3831.9575
$ws.Range("I132").Value = 2946.3901
$ws.Range("K132").Value = 8839.1703
$ws.Range("M132").Value = -6309.1703
